$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 845, shifting existing rows 845:886 down to 846:887
$ws.Rows.Item(845).Insert()

# Populate the newly inserted row 845 with the new data point.
# Force column A to be stored as text (not auto-parsed into a date serial),
# matching the original sheet's plain-text date cells, then drop the
# temporary number-format so the cell stays unstyled like its neighbours.
$ws.Cells.Item(845, 1).NumberFormat = "@"
$ws.Cells.Item(845, 1).Value = "2026/02/24"
$ws.Cells.Item(845, 1).ClearFormats()
$ws.Cells.Item(845, 2).Value = "火"
$ws.Cells.Item(845, 3).Value = 6
$ws.Cells.Item(845, 4).Value = 38
